# Wijzigingen en aanvullingen lineaire regressie
#
# 1. Rename the worksheet "economiedata" -> "data"
# 2. Remove the external workbook link (to ...\resources\grafieken.xlsx)
#    together with its cached external-link data (Edit Links > Break Link)
# 3. Move the active selection on the sheet to A28

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the sheet ---
$ws.Name = "data"

# --- 2. Break the external reference(s) so the link + cached data is gone ---
$linkSources = $wb.LinkSources(1)   # xlExcelLinks = 1
if ($linkSources) {
    foreach ($link in $linkSources) {
        $wb.BreakLink($link, 1) | Out-Null
    }
}

# --- 3. Update the selected cell on the sheet ---
$ws.Range("A28").Select() | Out-Null
